$wb = $excel.ActiveWorkbook

# Helper scratch cell used to force text (not auto-numeric) values onto
# target cells without disturbing the target cell's existing style index.
# (Range.Value with a purely-numeric-looking string auto-converts to a
# number; staging the text in a text-formatted scratch cell and copying
# only the *value* over via PasteSpecial keeps the destination's original
# formatting untouched.)
$scratchSheet = $wb.Worksheets.Item("Sheet1")
$scratch = $scratchSheet.Range("ZZ1000")
$scratch.NumberFormat = "@"

function Set-TextValue($range, $text) {
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163) # xlPasteValues
}

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet3 = $wb.Worksheets.Item("Sheet3")
$sheet4 = $wb.Worksheets.Item("Sheet4")

# Common change across all 4 sheets (row 2)
foreach ($ws in @($sheet1, $sheet2, $sheet3, $sheet4)) {
    Set-TextValue $ws.Range("F2") "5111367868"
    Set-TextValue $ws.Range("AE2") "8066212809"
    Set-TextValue $ws.Range("AT2") "9498816864"
    Set-TextValue $ws.Range("AX2") "8501155441"
}

# Sheet1-only additional changes
Set-TextValue $sheet1.Range("N2") "2024-03-06"
Set-TextValue $sheet1.Range("O2") "02:35:55 PM"
Set-TextValue $sheet1.Range("P2") "2024-03-06 07:57:38 PM"
Set-TextValue $sheet1.Range("AC2") "2024-03-06"
Set-TextValue $sheet1.Range("AN2") "97692"

# Clean up the scratch cell so it doesn't leave stray data behind.
$scratch.Clear()
